$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

$ws.Range("A1").Value = "Name: Germán"
$ws.Range("B1").Value = "Date: Mon Jan 01 2024 21:00:00 GMT-0300 (Argentina Standard Time)"

$ws.Range("A2").Value = "ID: 001"
$ws.Range("B2").Value = "Qty: 2"
$ws.Range("C2").Value = "Price: 120.5"
$ws.Range("D2").Value = "Item date: Tue Jan 02 2024 21:00:00 GMT-0300 (Argentina Standard Time)"
$ws.Range("E2").Value = "Missing: [[items.missingProp]]"
$ws.Range("F2").Formula = "=B2*C2"

$ws.Range("A3").Value = "ID: 002"
$ws.Range("B3").Value = "Qty: 0"
$ws.Range("C3").Value = "Price: "
$ws.Range("D3").Value = "Item date: "
$ws.Range("E3").Value = "Missing: [[items.missingProp]]"
$ws.Range("F3").Formula = "=B2*C2"
